$wb = $excel.ActiveWorkbook

# Worksheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 684
$ws.Range("F3").Value = 1515
$ws.Range("F4").Value = 3265
$ws.Range("F6").Value = 695
$ws.Range("F7").Value = 2276
$ws.Range("F8").Value = 494
$ws.Range("F9").Value = 412
$ws.Range("F12").Value = 336
$ws.Range("F14").Value = 444
$ws.Range("F15").Value = 15
$ws.Range("F17").Value = 231
$ws.Range("F18").Value = 4596
$ws.Range("F19").Value = 12
$ws.Range("F20").Value = 1316
$ws.Range("F21").Value = 3458
$ws.Range("F22").Value = 109
$ws.Range("F24").Value = 3626
$ws.Range("F25").Value = 5039
$ws.Range("F28").Value = 552
$ws.Range("F29").Value = 3244
$ws.Range("F30").Value = 366
$ws.Range("F34").Value = 880
$ws.Range("F35").Value = 1180
$ws.Range("F36").Value = 4
$ws.Range("F37").Value = 5
$ws.Range("F38").Value = 1418
$ws.Range("F39").Value = 126
$ws.Range("F40").Value = 1360
$ws.Range("F41").Value = 869
$ws.Range("F42").Value = 828
$ws.Range("F43").Value = 504
$ws.Range("F44").Value = 56
$ws.Range("F45").Value = 316
$ws.Range("F47").Value = 160
$ws.Range("F48").Value = 369
$ws.Range("F49").Value = 3723

# Worksheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 1006

# Worksheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2193

# Worksheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2193
$ws.Range("F3").Value = 684
$ws.Range("F4").Value = 1515
$ws.Range("F5").Value = 3265
$ws.Range("F7").Value = 695
$ws.Range("F9").Value = 2276
$ws.Range("F10").Value = 494
$ws.Range("F11").Value = 412
$ws.Range("F13").Value = 1006
$ws.Range("F16").Value = 336
$ws.Range("F18").Value = 444
$ws.Range("F19").Value = 15
$ws.Range("F20").Value = 231
$ws.Range("F21").Value = 4596
$ws.Range("F22").Value = 1316
$ws.Range("F24").Value = 3458
$ws.Range("F25").Value = 3626
$ws.Range("F26").Value = 5039
$ws.Range("F29").Value = 552
$ws.Range("F30").Value = 366
$ws.Range("F34").Value = 880
$ws.Range("F35").Value = 1180
$ws.Range("F37").Value = 1418
$ws.Range("F38").Value = 126
$ws.Range("F39").Value = 1360
$ws.Range("F40").Value = 869
$ws.Range("F43").Value = 56
$ws.Range("F45").Value = 316
$ws.Range("F47").Value = 160
$ws.Range("F48").Value = 369
$ws.Range("F49").Value = 3723
